# Insert a new weekly record for "Poroto granado" at row 13 (Fecha 2022-04-19),
# pushing the existing rows 13-92 down to 14-93.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = 8
$ws.Range("B13").Value = "Terminal La Palmera de La Serena"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44670
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 100112030
$ws.Range("G13").Value = "Poroto granado"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 600
$ws.Range("K13").Value = 26000
$ws.Range("L13").Value = 27000
$ws.Range("M13").Value = 26500
$ws.Range("N13").Value = '$/malla 25 kilos'
$ws.Range("O13").Value = "Provincia del Elquí"
$ws.Range("P13").Value = 1060
$ws.Range("Q13").Value = 25
$ws.Range("R13").Value = "Hortaliza"
